# Applies the profit-table recalculation updates across multiple sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) as produced by the scheduled
# Garuda_Profits runner. Each row's currentAveragePrice / LevePrice /
# LeveProfit columns (H-N) are refreshed to newly computed market values.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1944.7894
$ws.Range("I98").Value = 1830.6111
$ws.Range("J98").Value = 4000
$ws.Range("K98").Value = 1830.6111
$ws.Range("L98").Value = 4000
$ws.Range("M98").Value = -332.6111000000001
$ws.Range("N98").Value = -6996
$ws.Range("H122").Value = 1944.7894
$ws.Range("I122").Value = 1830.6111
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 5491.8333
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -3041.8333
$ws.Range("N122").Value = -16900
$ws.Range("H125").Value = 1279.5
$ws.Range("I125").Value = 625
$ws.Range("J125").Value = 1934
$ws.Range("K125").Value = 5625
$ws.Range("L125").Value = 17406
$ws.Range("M125").Value = -3165
$ws.Range("N125").Value = -22326
$ws.Range("H129").Value = 2753.415
$ws.Range("I129").Value = 467.08334
$ws.Range("J129").Value = 3422.5854
$ws.Range("K129").Value = 1401.25002
$ws.Range("L129").Value = 10267.7562
$ws.Range("M129").Value = 3598.74998
$ws.Range("N129").Value = -20267.7562

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1800
$ws.Range("I2").Value = 1500
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 1500
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = -1387
$ws.Range("N2").Value = -2226
$ws.Range("H32").Value = 13881.932
$ws.Range("I32").Value = 14107.262
$ws.Range("K32").Value = 14107.262
$ws.Range("M32").Value = -13820.262
$ws.Range("H45").Value = 795.375
$ws.Range("I45").Value = 724.8333
$ws.Range("J45").Value = 1007
$ws.Range("K45").Value = 724.8333
$ws.Range("L45").Value = 1007
$ws.Range("M45").Value = -347.8333
$ws.Range("N45").Value = -1761
$ws.Range("H110").Value = 1436.4375
$ws.Range("I110").Value = 950
$ws.Range("J110").Value = 2061.8572
$ws.Range("K110").Value = 950
$ws.Range("L110").Value = 2061.8572
$ws.Range("M110").Value = 1095
$ws.Range("N110").Value = -6151.8572
$ws.Range("H116").Value = 1800
$ws.Range("I116").Value = 1500
$ws.Range("J116").Value = 2000
$ws.Range("K116").Value = 1500
$ws.Range("L116").Value = 2000
$ws.Range("M116").Value = 794
$ws.Range("N116").Value = -6588
$ws.Range("H122").Value = 2136.8125
$ws.Range("I122").Value = 2108.0908
$ws.Range("K122").Value = 6324.2724
$ws.Range("M122").Value = -3874.2724

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1800
$ws.Range("I3").Value = 1500
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 1500
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = -1386
$ws.Range("N3").Value = -2228

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6062257
$ws.Range("I31").Value = 3128.9333
$ws.Range("J31").Value = 8334430
$ws.Range("K31").Value = 3128.9333
$ws.Range("L31").Value = 8334430
$ws.Range("M31").Value = -2833.9333
$ws.Range("N31").Value = -8335020
$ws.Range("H34").Value = 6062257
$ws.Range("I34").Value = 3128.9333
$ws.Range("J34").Value = 8334430
$ws.Range("K34").Value = 3128.9333
$ws.Range("L34").Value = 8334430
$ws.Range("M34").Value = -2926.9333
$ws.Range("N34").Value = -8334834
$ws.Range("H127").Value = 38700
$ws.Range("J127").Value = 38700
$ws.Range("L127").Value = 38700
$ws.Range("N127").Value = -48620
$ws.Range("H134").Value = 861.1070999999999
$ws.Range("I134").Value = 864.4400000000001
$ws.Range("K134").Value = 2593.32
$ws.Range("M134").Value = -58.32000000000016

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 596.6875
$ws.Range("I107").Value = 367.77777
$ws.Range("J107").Value = 891
$ws.Range("K107").Value = 1103.33331
$ws.Range("L107").Value = 2673
$ws.Range("M107").Value = 816.66669
$ws.Range("N107").Value = -6513
$ws.Range("H122").Value = 397720.3
$ws.Range("I122").Value = 392.64706
$ws.Range("J122").Value = 1011772.2
$ws.Range("K122").Value = 3533.82354
$ws.Range("L122").Value = 9105949.799999999
$ws.Range("M122").Value = -1083.82354
$ws.Range("N122").Value = -9110849.799999999
$ws.Range("H131").Value = 3741342
$ws.Range("I131").Value = 11503.333
$ws.Range("J131").Value = 5339844.5
$ws.Range("K131").Value = 34509.999
$ws.Range("L131").Value = 16019533.5
$ws.Range("M131").Value = -29469.999
$ws.Range("N131").Value = -16029613.5
$ws.Range("H132").Value = 2194.2666
$ws.Range("J132").Value = 2628.182
$ws.Range("L132").Value = 23653.638
$ws.Range("N132").Value = -28713.638

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 25003520
$ws.Range("I113").Value = 41671060
$ws.Range("J113").Value = 2210
$ws.Range("K113").Value = 41671060
$ws.Range("L113").Value = 2210
$ws.Range("M113").Value = -41668890
$ws.Range("N113").Value = -6550
$ws.Range("H122").Value = 6299.933
$ws.Range("I122").Value = 15424.75
$ws.Range("J122").Value = 2981.818
$ws.Range("K122").Value = 46274.25
$ws.Range("L122").Value = 8945.454000000002
$ws.Range("M122").Value = -43824.25
$ws.Range("N122").Value = -13845.454

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2000
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("H36").Value = 61900
$ws.Range("J36").Value = 61900
$ws.Range("L36").Value = 61900
$ws.Range("N36").Value = -63024
$ws.Range("H61").Value = 3152
$ws.Range("I61").Value = 3152
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3152
$ws.Range("L61").Value = 0
$ws.Range("N61").Value = -2950
$ws.Range("H113").Value = 3152
$ws.Range("I113").Value = 3152
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3152
$ws.Range("L113").Value = 0
$ws.Range("N113").Value = -982
$ws.Range("H122").Value = 4986.2856
$ws.Range("I122").Value = 17004
$ws.Range("J122").Value = 2983.3333
$ws.Range("K122").Value = 51012
$ws.Range("L122").Value = 8949.999899999999
$ws.Range("M122").Value = -48562
$ws.Range("N122").Value = -13849.9999
$ws.Range("H126").Value = 2000
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("M61").ClearContents()
$ws.Range("M113").ClearContents()
$ws.Range("N126").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1000
$ws.Range("I113").Value = 1300
$ws.Range("J113").Value = 760
$ws.Range("K113").Value = 3900
$ws.Range("L113").Value = 2280
$ws.Range("M113").Value = -1730
$ws.Range("N113").Value = -6620
$ws.Range("H126").Value = 3650.3333
$ws.Range("I126").Value = 4250.5
$ws.Range("J126").Value = 2450
$ws.Range("K126").Value = 12751.5
$ws.Range("L126").Value = 7350
$ws.Range("M126").Value = -10281.5
$ws.Range("N126").Value = -12290

